$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.551048874855042
$ws.Range("B1").Value = -1
$ws.Range("C1").Value = 4.572749614715576
$ws.Range("D1").Value = 2.273092985153198
$ws.Range("E1").Value = 1.808940649032593
